# Applies a weekly reshuffle of the Fecha/Volumen/Precio columns (D, M, N, O, P, S)
# across data rows 2-30 of the sheet, per the commit "Fruta / hortaliza, semanal".
# Columns A,B,C,E,F,G,H,I,J,K,L,Q,R,T are left untouched; only D,M,N,O,P,S move
# between rows according to the mapping below (target row -> source row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# target row number -> source row number (where the new data used to live)
$rowMap = @{
    2  = 29
    3  = 17
    4  = 5
    5  = 4
    6  = 22
    7  = 19
    8  = 23
    9  = 25
    10 = 8
    11 = 3
    12 = 16
    13 = 6
    14 = 11
    15 = 10
    16 = 14
    17 = 15
    18 = 18
    19 = 21
    20 = 9
    21 = 26
    22 = 27
    23 = 2
    24 = 7
    25 = 13
    26 = 30
    27 = 28
    28 = 20
    29 = 24
    30 = 12
}

$cols = @("D", "M", "N", "O", "P", "S")

# Snapshot original values for the affected columns before overwriting anything.
$snapshot = @{}
for ($r = 2; $r -le 30; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# Write the shuffled values back out.
foreach ($targetRow in $rowMap.Keys) {
    $sourceRow = $rowMap[$targetRow]
    $srcVals = $snapshot[$sourceRow]
    foreach ($c in $cols) {
        $ws.Range("$c$targetRow").Value = $srcVals[$c]
    }
}
